$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111902031
$ws.Range("B2").Value = 90808
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2"
$ws.Range("Q2").Value = 524990
$ws.Range("R2").Value = 6867386
# Row 3
$ws.Range("A3").Value = 111902037
$ws.Range("B3").Value = 90802
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 149
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "Tallgråticka"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "Boletopsis grisea"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "(Peck) Bondartsev & Singer"
$ws.Range("Q3").Value = 524869
$ws.Range("R3").Value = 6867441
$ws.Range("S3").Value = 5
# Row 4
$ws.Range("A4").Value = 111902040
$ws.Range("B4").Value = 90448
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 4745
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "Tallriska"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "Lactarius musteus"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "Fr."
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "1"
$ws.Range("Q4").Value = 524891
$ws.Range("R4").Value = 6866840
$ws.Range("S4").Value = 10
# Row 5
$ws.Range("A5").Value = 111902036
$ws.Range("B5").Value = 88180
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 6276
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "Goliatmusseron"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "Tricholoma matsutake"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "2"
$ws.Range("Q5").Value = 525016
$ws.Range("R5").Value = 6867406
$ws.Range("S5").Value = 25
# Row 6
$ws.Range("A6").Value = 111902038
$ws.Range("B6").Value = 90814
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 4364
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "Dropptaggsvamp"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "Hydnellum ferrugineum"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "1"
$ws.Range("Q6").Value = 524893
$ws.Range("R6").Value = 6867499
$ws.Range("S6").Value = 10
# Row 7
$ws.Range("A7").Value = 111902033
$ws.Range("B7").Value = 90448
$ws.Range("E7").Value = 4745
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "Tallriska"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "Lactarius musteus"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "Fr."
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "1"
$ws.Range("Q7").Value = 525027
$ws.Range("R7").Value = 6867370
$ws.Range("S7").Value = 10
# Row 8
$ws.Range("A8").Value = 111902030
$ws.Range("B8").Value = 88180
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 6276
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "Goliatmusseron"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "Tricholoma matsutake"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "6"
$ws.Range("Q8").Value = 524971
$ws.Range("R8").Value = 6867379
$ws.Range("S8").Value = 5
# Row 9
$ws.Range("A9").Value = 111902027
$ws.Range("B9").Value = 90808
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 4362
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "Blå taggsvamp"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "Hydnellum caeruleum"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "(Hornem.) P.Karst."
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "5"
$ws.Range("Q9").Value = 524937
$ws.Range("R9").Value = 6867322
$ws.Range("S9").Value = 25
# Row 10
$ws.Range("A10").Value = 111902026
$ws.Range("B10").Value = 90830
$ws.Range("E10").Value = 2059
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "Skrovlig taggsvamp"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "Hydnellum scabrosum"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = "1"
$ws.Range("Q10").Value = 524951
$ws.Range("R10").Value = 6867324
$ws.Range("S10").Value = 10
# Row 11
$ws.Range("A11").Value = 111902035
$ws.Range("B11").Value = 90806
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 4361
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "Orange taggsvamp"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "Hydnellum aurantiacum"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "3"
$ws.Range("Q11").Value = 525047
$ws.Range("R11").Value = 6867385
$ws.Range("S11").Value = 25
# Row 12
$ws.Range("A12").Value = 111902039
$ws.Range("B12").Value = 90830
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 2059
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "Skrovlig taggsvamp"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "Hydnellum scabrosum"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I12").NumberFormat = "@"
$ws.Range("I12").Value = "5"
$ws.Range("Q12").Value = 524868
$ws.Range("R12").Value = 6867460
# Row 13
$ws.Range("A13").Value = 111902028
$ws.Range("B13").Value = 90814
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 4364
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "Dropptaggsvamp"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "Hydnellum ferrugineum"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = "1"
$ws.Range("Q13").Value = 524954
$ws.Range("R13").Value = 6867304
$ws.Range("S13").Value = 5
# Row 14
$ws.Range("A14").Value = 111902029
$ws.Range("B14").Value = 88180
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "VU"
$ws.Range("E14").Value = 6276
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "Goliatmusseron"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "Tricholoma matsutake"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = "4"
$ws.Range("Q14").Value = 524972
$ws.Range("R14").Value = 6867342
$ws.Range("S14").Value = 5
# Row 15
$ws.Range("B15").Value = 90806
# Row 16
$ws.Range("A16").Value = 111902034
$ws.Range("B16").Value = 90808
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 4362
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "Blå taggsvamp"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "Hydnellum caeruleum"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "(Hornem.) P.Karst."
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "10"
$ws.Range("Q16").Value = 525039
$ws.Range("R16").Value = 6867407
$ws.Range("S16").Value = 25
